$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee")

# Update Username column (column C) values: append "5" to each username
$ws.Range("C2").Value = "msamuels1235"
$ws.Range("C3").Value = "sjones125"
$ws.Range("C4").Value = "wjohnson1235"
$ws.Range("C5").Value = "asmith125"
$ws.Range("C6").Value = "jshmit1235"

# Move the active cell selection to G12
$ws.Range("G12").Select()
